# Applies the "added the GUI component" commit's data changes to the
# trading_synthesis workbook: refreshed indicator values on the weekly,
# daily and 4h sheets, a new note on weekly!J4, and six freshly logged
# 4h signals (rows 10-14, with the old row 9 pushed down and restated).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "weekly"
# ---------------------------------------------------------------------
$weekly = $wb.Worksheets.Item("weekly")

$weekly.Range("E2").Value = 123.64
$weekly.Range("F2").Value = 36.52
$weekly.Range("I2").Value = 0

$weekly.Range("D3").Value = 193.5
$weekly.Range("E3").Value = 101.94
$weekly.Range("F3").Value = 63.25

$weekly.Range("D4").Value = 1372.45
$weekly.Range("E4").Value = -173.54
$weekly.Range("F4").Value = 79.41
$weekly.Range("J4").Value = "test"

$weekly.Range("D8").Value = 107702.3
$weekly.Range("E8").Value = 127.18
$weekly.Range("F8").Value = 79.72

$weekly.Range("D9").Value = 76.76000000000001
$weekly.Range("E9").Value = 139.9
$weekly.Range("F9").Value = 69.95

$weekly.Range("D10").Value = 244.23
$weekly.Range("E10").Value = 219.96
$weekly.Range("F10").Value = 87.09

$weekly.Range("D11").Value = 23.27
$weekly.Range("E11").Value = -117.47

$weekly.Range("D12").Value = 1208.15
$weekly.Range("E12").Value = 148.84
$weekly.Range("F12").Value = 90.87

$weekly.Range("D13").Value = 459.83
$weekly.Range("E13").Value = 143.35
$weekly.Range("F13").Value = 63.16

$weekly.Range("D14").Value = 253.36
$weekly.Range("E14").Value = 184.03
$weekly.Range("F14").Value = 49.16
$weekly.Range("G14").Value = 52.29

$weekly.Range("D15").Value = 277.42
$weekly.Range("E15").Value = 114.83
$weekly.Range("F15").Value = 62.22

$weekly.Range("D16").Value = 1.53
$weekly.Range("E16").Value = 112.28
$weekly.Range("F16").Value = 8.25
$weekly.Range("H16").Value = -0.06

# ---------------------------------------------------------------------
# Sheet "daily"
# ---------------------------------------------------------------------
$daily = $wb.Worksheets.Item("daily")

$daily.Range("D2").Value = 8744.559999999999
$daily.Range("E2").Value = 107.32
$daily.Range("F2").Value = 63.21
$daily.Range("G2").Value = 67.06

$daily.Range("A4").Value = 45805
$daily.Range("D4").Value = 108.57
$daily.Range("E4").Value = 110.81
$daily.Range("F4").Value = 52.27
$daily.Range("G4").Value = 67.69
$daily.Range("H4").Value = -0.77
$daily.Range("I4").Value = -0.63

$daily.Range("A5").Value = 45805
$daily.Range("D5").Value = 4.09
$daily.Range("E5").Value = 160.79
$daily.Range("F5").Value = 16.2
$daily.Range("G5").Value = 21.68
$daily.Range("H5").Value = 0.55
$daily.Range("I5").Value = -0.92

$daily.Range("A6").Value = 45805
$daily.Range("D6").Value = 9.83
$daily.Range("E6").Value = 168.8
$daily.Range("F6").Value = 32.06
$daily.Range("G6").Value = 34.6
$daily.Range("H6").Value = -0.19

# ---------------------------------------------------------------------
# Sheet "4h"
# ---------------------------------------------------------------------
$h4 = $wb.Worksheets.Item("4h")

$h4.Range("E2").Value = -111.13
$h4.Range("F2").Value = 51.31

$h4.Range("E3").Value = -104.44
$h4.Range("F3").Value = 70.48

$h4.Range("B4").Value = "Sell"
$h4.Range("C4").Value = "EURJPY=X"
$h4.Range("D4").Value = 163.67
$h4.Range("E4").Value = 101.23
$h4.Range("F4").Value = 34.74
$h4.Range("G4").Value = 42.61
$h4.Range("H4").Value = 0.21
$h4.Range("I4").Value = -0.63

$h4.Range("A5").Value = 45805.5
$h4.Range("C5").Value = "GBPUSD=X"
$h4.Range("D5").Value = 1.35
$h4.Range("E5").Value = -115.83
$h4.Range("F5").Value = 81.53
$h4.Range("G5").Value = 67.97
$h4.Range("H5").Value = 0.3
$h4.Range("I5").Value = 0.99

$h4.Range("A6").Value = 45805.5
$h4.Range("B6").Value = "Sell"
$h4.Range("C6").Value = "USDJPY=X"
$h4.Range("D6").Value = 144.57
$h4.Range("E6").Value = 107.88
$h4.Range("F6").Value = 27.72
$h4.Range("G6").Value = 43.81
$h4.Range("H6").Value = -0.31
$h4.Range("I6").Value = -1

$h4.Range("A7").Value = 45805.5
$h4.Range("C7").Value = "USDILS=X"
$h4.Range("D7").Value = 3.51
$h4.Range("E7").Value = -108.46
$h4.Range("F7").Value = 36.85
$h4.Range("G7").Value = 28.97
$h4.Range("H7").Value = -0.05
$h4.Range("I7").Value = 0.51

$h4.Range("A8").Value = 45805.54166666666
$h4.Range("B8").Value = "Buy"
$h4.Range("C8").Value = "^SSMI"
$h4.Range("D8").Value = 12199.01
$h4.Range("E8").Value = -108.22
$h4.Range("F8").Value = 81.56
$h4.Range("G8").Value = 77.36
$h4.Range("H8").Value = -0.07000000000000001
$h4.Range("I8").Value = 0.67

$h4.Range("A9").Value = 45805.39583333334
$h4.Range("C9").Value = "LIT"
$h4.Range("D9").Value = 37.08
$h4.Range("E9").Value = -132.36
$h4.Range("F9").Value = 64.90000000000001
$h4.Range("G9").Value = 49.1
$h4.Range("H9").Value = 0.36
$h4.Range("I9").Value = 0.97

# New rows 10-14, freshly logged 4h signals.
$h4.Range("A10").Value = 45805.39583333334
$h4.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$h4.Range("B10").Value = "Buy"
$h4.Range("C10").Value = "GM"
$h4.Range("D10").Value = 48.47
$h4.Range("E10").Value = -149.89
$h4.Range("F10").Value = 62.58
$h4.Range("G10").Value = 43.7
$h4.Range("H10").Value = 0.96
$h4.Range("I10").Value = 0.57

$h4.Range("A11").Value = 45805.39583333334
$h4.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$h4.Range("B11").Value = "Buy"
$h4.Range("C11").Value = "LYV"
$h4.Range("D11").Value = 141.87
$h4.Range("E11").Value = -119.52
$h4.Range("F11").Value = 87.22
$h4.Range("G11").Value = 73.89
$h4.Range("H11").Value = 0.32
$h4.Range("I11").Value = 1

$h4.Range("A12").Value = 45805.39583333334
$h4.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$h4.Range("B12").Value = "Buy"
$h4.Range("C12").Value = "RIOT"
$h4.Range("D12").Value = 8.48
$h4.Range("E12").Value = -107.21
$h4.Range("F12").Value = 67.13
$h4.Range("G12").Value = 46.36
$h4.Range("H12").Value = 1.07
$h4.Range("I12").Value = 1.11

$h4.Range("A13").Value = 45805.39583333334
$h4.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$h4.Range("B13").Value = "Buy"
$h4.Range("C13").Value = "BITF"
$h4.Range("D13").Value = 1.01
$h4.Range("E13").Value = -172.7
$h4.Range("F13").Value = 75.61
$h4.Range("G13").Value = 55.51
$h4.Range("H13").Value = 0.9
$h4.Range("I13").Value = 1.29

$h4.Range("A14").Value = 45805.39583333334
$h4.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$h4.Range("B14").Value = "Buy"
$h4.Range("C14").Value = "ARBK"
$h4.Range("D14").Value = 0.4
$h4.Range("E14").Value = -104.58
$h4.Range("F14").Value = 66.52
$h4.Range("G14").Value = 51.41
$h4.Range("H14").Value = 0.6899999999999999
$h4.Range("I14").Value = 0.9399999999999999
